$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: R_I2C2/R_I2C1/... resistor group ---
# Designators: drop R_SET1,R_SET4,R_SET3,R_SET2 and fold in R_PROG1 (its old,
# now-removed standalone row below used the same 10k Ohm part).
$ws.Range("B3").Value = "R_I2C2,R_I2C1,R_ADDR1,R_PROG1"
$ws.Range("E3").Value = "10k Ohm"
$ws.Range("F3").Value = "TE Connectivity / Holsworthy "
$ws.Range("G3").Value = "CRG0603J10K "

# --- Row 4: C2,C1 100nF capacitor - remove manufacturer / part number info ---
$ws.Range("F4:G4").ClearContents()

# --- Remove the now-redundant standalone R_PROG1 row (old row 7) ---
# Move J_USB1 (old row 8) up into row 7, then clear out the vacated row 8.
for ($c = 1; $c -le 7; $c++) {
    $ws.Cells.Item(7, $c).Value = $ws.Cells.Item(8, $c).Value2
}
$ws.Rows("8:8").Clear()
$ws.Rows("7:7").RowHeight = 15

# --- Clear the stray empty formatted cells left in row 9 ---
$ws.Range("B9:E9").Clear()

# --- Normalize the highlighted manufacturer/part cells to the common centered
# + wrapped body style (matches the rest of the table) instead of the old
# center/no-wrap variant. Copy format only, so no new cell styles are created.
$normalized = @("F1","G2","F3","G3","F4","G5","G6","G7")
foreach ($addr in $normalized) {
    $ws.Range("E1").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# --- Restore the active cell/selection to C3 ---
$null = $ws.Range("C3").Select()
